# dict_detail.xlsx template: add 4 new "aggregate/foreign-tab" columns
# (create_usr_id, create_time, update_usr_id, update_time) to both the
# field-comment header row (row 1, driven by `comment.*`) and the
# data-model sample row (row 2, driven by `model.*`).
#
# Columns A-G are left untouched; H-K are appended after them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - header/comment placeholders (with inline <%_dataValidation_(...)%> setup
# for the two "usr_id" select-list columns, mirroring the existing is_locked /
# is_enabled / dict_id columns but without the allowBlank override).
$ws.Range("H1").Value = '<%=comment.create_usr_id_lbl%><%selectList.create_usr_id = data.findAllUsr.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.create_usr_id.join(",") }"` })%>'
$ws.Range("I1").Value = '<%=comment.create_time_lbl%>'
$ws.Range("J1").Value = '<%=comment.update_usr_id_lbl%><%selectList.update_usr_id = data.findAllUsr.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.update_usr_id.join(",") }"` })%>'
$ws.Range("K1").Value = '<%=comment.update_time_lbl%>'

# Row 2 - sample/data-model row placeholders.
$ws.Range("H2").Value = '<%=model.create_usr_id_lbl%>'
$ws.Range("I2").Value = '<%~model.create_time ? new Date(model.create_time) : ""%>'
$ws.Range("J2").Value = '<%=model.update_usr_id_lbl%>'
$ws.Range("K2").Value = '<%~model.update_time ? new Date(model.update_time) : ""%>'
